# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.640.11'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '2.655.97'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'596.22"
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').Value = "'158.41"
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('D7').Value = "'0.642"
$ws.Range('E7').Value = '  +4.65%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('D10').Value = "'5.83"
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('D11').Value = "'0.397"
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').Value = "'28.99"
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').Value = '3.131.41'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').Value = '65.521.76'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '2.670.05'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = "'12.46"
$ws.Range('E18').Value = '  -3.34%  '
$ws.Range('D19').Value = "'4.76"
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').Value = "'352.55"
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = "'7.43"
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('D22').Value = "'1.00"
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = "'69.41"
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D24').Value = "'0.0000113"
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').Value = "'1.76"
$ws.Range('E25').Value = '  +4.06%  '
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').Value = "'1.61"
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('D28').Value = "'562.02"
$ws.Range('E28').Value = '  +5.77%  '
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('D30').Value = "'8.06"
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = "'0.995"
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('D33').Value = "'1.80"
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('D34').Value = "'6.66"
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('D35').Value = "'5.45"
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('D36').Value = "'0.421"
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').Value = "'20.45"
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('D38').Value = "'0.998"
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = "'1.96"
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('D40').Value = "'152.88"
$ws.Range('E40').Value = '  -3.30%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = "'160.99"
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = "'2.45"
$ws.Range('E43').Value = '  +4.43%  '
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').Value = "'0.0613"
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').Value = "'23.19"
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').Value = "'19.65"
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('D51').Value = '0.0₆0244'
$ws.Range('E51').Value = '  -7.74%  '
